$wb = $excel.ActiveWorkbook

# Sheet "展览" (1st sheet): update 想去人数 (want-to-go count) for two events
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F4").Value = 1157
$wsExhibit.Range("F5").Value = 594

# Sheet "全部类型" (4th sheet): same two events appear again, update accordingly
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 1157
$wsAll.Range("F6").Value = 594
